# Auto-generated Excel COM-interop script applying the Leve profit-value refresh
# from the scheduled price-data runner (see commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 35157292
$ws.Range("I19").Value = 59210950
$ws.Range("J19").Value = 1943.8462
$ws.Range("K19").Value = 59210950
$ws.Range("L19").Value = 1943.8462
$ws.Range("M19").Value = -59210775
$ws.Range("N19").Value = -2293.8462

# Row 32
$ws.Range("H32").Value = 888.26666
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 993.0909
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 993.0909
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1645.0909

# Row 137
$ws.Range("H137").Value = 37991692
$ws.Range("I137").Value = 8334461.5
$ws.Range("J137").Value = 186277840
$ws.Range("K137").Value = 25003384.5
$ws.Range("L137").Value = 558833520
$ws.Range("M137").Value = -25000834.5
$ws.Range("N137").Value = -558838620

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 1512
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 32
$ws.Range("H32").Value = 2459969.8
$ws.Range("I32").Value = 3213878.8
$ws.Range("K32").Value = 3213878.8
$ws.Range("M32").Value = -3213591.8

# Row 61
$ws.Range("H61").Value = 3590505
$ws.Range("I61").Value = 3206140.2
$ws.Range("J61").Value = 3923621.2
$ws.Range("K61").Value = 3206140.2
$ws.Range("L61").Value = 3923621.2
$ws.Range("M61").Value = -3205928.2
$ws.Range("N61").Value = -3924045.2

# Row 74
$ws.Range("H74").Value = 48573880
$ws.Range("I74").Value = 63394000
$ws.Range("J74").Value = 22227004
$ws.Range("K74").Value = 63394000
$ws.Range("L74").Value = 22227004
$ws.Range("M74").Value = -63393126
$ws.Range("N74").Value = -22228752

# Row 77
$ws.Range("H77").Value = 48573880
$ws.Range("I77").Value = 63394000
$ws.Range("J77").Value = 22227004
$ws.Range("K77").Value = 316970000
$ws.Range("L77").Value = 111135020
$ws.Range("M77").Value = -316965632
$ws.Range("N77").Value = -111143756

# Row 132
$ws.Range("H132").Value = 13414852
$ws.Range("I132").Value = 13338114
$ws.Range("J132").Value = 13894464
$ws.Range("K132").Value = 40014342
$ws.Range("L132").Value = 41683392
$ws.Range("M132").Value = -40011812
$ws.Range("N132").Value = -41688452

# Row 136
$ws.Range("H136").Value = 3590505
$ws.Range("I136").Value = 3206140.2
$ws.Range("J136").Value = 3923621.2
$ws.Range("K136").Value = 9618420.600000001
$ws.Range("L136").Value = 11770863.6
$ws.Range("M136").Value = -9615870.600000001
$ws.Range("N136").Value = -11775963.6

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1083.3334
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -5090

# Row 134
$ws.Range("H134").Value = 15307198
$ws.Range("I134").Value = 18519364
$ws.Range("J134").Value = 4466136.5
$ws.Range("K134").Value = 55558092
$ws.Range("L134").Value = 13398409.5
$ws.Range("M134").Value = -55555557
$ws.Range("N134").Value = -13403479.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5212923
$ws.Range("J31").Value = 12513111
$ws.Range("L31").Value = 12513111
$ws.Range("N31").Value = -12513701

# Row 34
$ws.Range("H34").Value = 5212923
$ws.Range("J34").Value = 12513111
$ws.Range("L34").Value = 12513111
$ws.Range("N34").Value = -12513515

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3264021.8
$ws.Range("I5").Value = 2564550
$ws.Range("J5").Value = 4762890
$ws.Range("K5").Value = 7693650
$ws.Range("L5").Value = 14288670
$ws.Range("M5").Value = -7693538
$ws.Range("N5").Value = -14288894

# Row 12
$ws.Range("H12").Value = 97.77778000000001
$ws.Range("I12").Value = 108.57143
$ws.Range("J12").Value = 94
$ws.Range("K12").Value = 325.71429
$ws.Range("L12").Value = 282
$ws.Range("M12").Value = -152.71429
$ws.Range("N12").Value = -628

# Row 23
$ws.Range("H23").Value = 42.77778
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 49.285713
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 147.857139
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -617.857139

# Row 70
$ws.Range("H70").Value = 1290
$ws.Range("I70").Value = 1290
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3870
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3555
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 1290
$ws.Range("I73").Value = 1290
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3870
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2778
$ws.Range("N73").ClearContents()

# Row 135
$ws.Range("H135").Value = 3264021.8
$ws.Range("I135").Value = 2564550
$ws.Range("J135").Value = 4762890
$ws.Range("K135").Value = 23080950
$ws.Range("L135").Value = 42866010
$ws.Range("M135").Value = -23078415
$ws.Range("N135").Value = -42871080

# Row 141
$ws.Range("H141").Value = 2339.6428
$ws.Range("I141").Value = 2062.9167
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 6188.750100000001
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -1008.750100000001
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 1787
$ws.Range("I36").Value = 1800
$ws.Range("K36").Value = 1800
$ws.Range("M36").Value = -1315

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 588
$ws.Range("I10").Value = 486.66666
$ws.Range("J10").Value = 740
$ws.Range("K10").Value = 486.66666
$ws.Range("L10").Value = 740
$ws.Range("M10").Value = -346.66666
$ws.Range("N10").Value = -1020

# Row 17
$ws.Range("H17").Value = 1300004.5
$ws.Range("J17").Value = 66672.664
$ws.Range("L17").Value = 66672.664
$ws.Range("N17").Value = -67012.664

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 9252
$ws.Range("I8").Value = 12332.667
$ws.Range("J8").Value = 10
$ws.Range("K8").Value = 12332.667
$ws.Range("L8").Value = 10
$ws.Range("M8").Value = -12192.667
$ws.Range("N8").Value = -290

# Row 123
$ws.Range("H123").Value = 46452.25
$ws.Range("J123").Value = 46452.25
$ws.Range("L123").Value = 46452.25
$ws.Range("N123").Value = -56252.25

# Row 132
$ws.Range("H132").Value = 2494174.2
$ws.Range("I132").Value = 2030447.2
$ws.Range("K132").Value = 6091341.6
$ws.Range("M132").Value = -6088811.6
